$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add row 5 ("Crumpet GEF" / "Crumpet exporter") ---
# Copy formatting from row 2 (an existing fully-populated data row) down to row 5
# so the new row picks up the same cell styles used by the other data rows.
$ws.Range("A2:J2").Copy($ws.Range("A5:J5"))

$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# --- Fill row 6 ("Scone GEF" / "Scone exporter") ---
# Row 6 already exists (previously blank placeholder cells with the right
# styles already applied), so we only need to populate the values.
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# --- Update the active selection to reflect the newly added rows ---
$ws.Activate()
$null = $ws.Range("A5:J6").Select()
